$wb = $excel.ActiveWorkbook

# --- Rename sheets (tab names) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-1650291165529401"
$wb.Worksheets.Item(2).Name = "NB_TO-16502911669924002"
$wb.Worksheets.Item(3).Name = "RS_TO-16502911669954"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502911670543966"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650291167115406"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502911654724011.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911654953978.csv"
$ws1.Range("B4").Value = "go_stims-16502911654964354.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911655274012.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16502911669714012.csv"
$ws2.Range("B3").Value = "ZB-match_6-16502911657913961.csv"
$ws2.Range("B4").Value = "OB-1650291166536398.csv"
$ws2.Range("B5").Value = "TB-16502911667074335.csv"
$ws2.Range("B6").Value = "ZB-match_4-16502911657573974.csv"
$ws2.Range("B7").Value = "TB-1650291166939404.csv"
$ws2.Range("B8").Value = "OB-16502911665953987.csv"
$ws2.Range("B9").Value = "OB-16502911664724014.csv"
$ws2.Range("B10").Value = "ZB-match_5-16502911655624046.csv"

# --- Sheet 3 (RS) - swap eyes open / eyes closed ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502911670233994.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911669973984.csv"
$ws4.Range("B4").Value = "MM_stims-16502911670383992.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911670264156.csv"
$ws4.Range("B6").Value = "MM_stims-16502911670533986.csv"
$ws4.Range("B7").Value = "ZM_stims-1650291167039399.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1650291167057397.csv"
$ws5.Range("B3").Value = "SAT_stims-16502911670704048.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502911671003985.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502911670844.csv"
